$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1682539682539683
$ws.Range("C2").Value = 0.6063492063492063
$ws.Range("J2").Value = 0.01587301587301587
$ws.Range("P2").Value = 0.1047619047619048
$ws.Range("S2").Value = 0.1047619047619048
$ws.Range("B3").Value = 0.009852216748768473
$ws.Range("C3").Value = 0.02463054187192118
$ws.Range("J3").Value = 0.04433497536945813
$ws.Range("P3").Value = 0.7438423645320197
$ws.Range("S3").Value = 0.1773399014778325
$ws.Range("J4").Value = 0.08333333333333333
$ws.Range("P4").Value = 0.6388888888888888
$ws.Range("S4").Value = 0.2777777777777778
$ws.Range("B6").Value = 0.08365019011406843
$ws.Range("D6").Value = 0.007604562737642586
$ws.Range("F6").Value = 0.08745247148288973
$ws.Range("J6").Value = 0.2053231939163498
$ws.Range("O6").Value = 0.02281368821292776
$ws.Range("Q6").Value = 0.1368821292775665
$ws.Range("R6").Value = 0.09505703422053231
$ws.Range("S6").Value = 0.3612167300380228
$ws.Range("B7").Value = 0.1184210526315789
$ws.Range("D7").Value = 0.01973684210526316
$ws.Range("F7").Value = 0.05921052631578947
$ws.Range("J7").Value = 0.1052631578947368
$ws.Range("O7").Value = 0.01973684210526316
$ws.Range("Q7").Value = 0.1578947368421053
$ws.Range("R7").Value = 0.1052631578947368
$ws.Range("S7").Value = 0.4144736842105263
$ws.Range("B8").Value = 0.1138392857142857
$ws.Range("D8").Value = 0.02232142857142857
$ws.Range("E8").Value = 0.002232142857142857
$ws.Range("F8").Value = 0.07589285714285714
$ws.Range("J8").Value = 0.1160714285714286
$ws.Range("O8").Value = 0.006696428571428571
$ws.Range("Q8").Value = 0.1830357142857143
$ws.Range("R8").Value = 0.07142857142857142
$ws.Range("S8").Value = 0.4084821428571428
$ws.Range("B9").Value = 0.119815668202765
$ws.Range("D9").Value = 0.02304147465437788
$ws.Range("F9").Value = 0.1152073732718894
$ws.Range("J9").Value = 0.09216589861751152
$ws.Range("O9").Value = 0.009216589861751152
$ws.Range("Q9").Value = 0.1981566820276498
$ws.Range("R9").Value = 0.08294930875576037
$ws.Range("S9").Value = 0.3594470046082949
$ws.Range("B10").Value = 0.1096214511041009
$ws.Range("D10").Value = 0.01577287066246057
$ws.Range("E10").Value = 0.0007886435331230284
$ws.Range("F10").Value = 0.07728706624605679
$ws.Range("J10").Value = 0.1135646687697161
$ws.Range("O10").Value = 0.01419558359621451
$ws.Range("Q10").Value = 0.194794952681388
$ws.Range("R10").Value = 0.08675078864353312
$ws.Range("S10").Value = 0.3872239747634069
$ws.Range("G11").Value = 0.1399317406143345
$ws.Range("J11").Value = 0.08873720136518772
$ws.Range("K11").Value = 0.2116040955631399
$ws.Range("L11").Value = 0.5221843003412969
$ws.Range("S11").Value = 0.03754266211604096
$ws.Range("G12").Value = 0.6049382716049383
$ws.Range("J12").Value = 0.2654320987654321
$ws.Range("K12").Value = 0.02469135802469136
$ws.Range("L12").Value = 0.0308641975308642
$ws.Range("S12").Value = 0.07407407407407407
$ws.Range("G13").Value = 0.5882352941176471
$ws.Range("J13").Value = 0.2941176470588235
$ws.Range("S13").Value = 0.1176470588235294
$ws.Range("F15").Value = 0.01951219512195122
$ws.Range("H15").Value = 0.1268292682926829
$ws.Range("I15").Value = 0.06341463414634146
$ws.Range("J15").Value = 0.4682926829268293
$ws.Range("K15").Value = 0.08292682926829269
$ws.Range("M15").Value = 0.01463414634146342
$ws.Range("O15").Value = 0.03902439024390244
$ws.Range("S15").Value = 0.1853658536585366
$ws.Range("F16").Value = 0.005076142131979695
$ws.Range("H16").Value = 0.1522842639593909
$ws.Range("I16").Value = 0.07106598984771574
$ws.Range("J16").Value = 0.4060913705583756
$ws.Range("K16").Value = 0.1421319796954315
$ws.Range("M16").Value = 0.02538071065989848
$ws.Range("O16").Value = 0.06091370558375635
$ws.Range("S16").Value = 0.1370558375634518
$ws.Range("F17").Value = 0.0187793427230047
$ws.Range("H17").Value = 0.2230046948356808
$ws.Range("I17").Value = 0.1173708920187793
$ws.Range("J17").Value = 0.431924882629108
$ws.Range("K17").Value = 0.06572769953051644
$ws.Range("M17").Value = 0.007042253521126761
$ws.Range("O17").Value = 0.05164319248826291
$ws.Range("S17").Value = 0.08450704225352113
$ws.Range("F18").Value = 0.02487562189054726
$ws.Range("H18").Value = 0.1741293532338309
$ws.Range("I18").Value = 0.09950248756218906
$ws.Range("J18").Value = 0.3631840796019901
$ws.Range("K18").Value = 0.1243781094527363
$ws.Range("M18").Value = 0.004975124378109453
$ws.Range("N18").Value = 0.004975124378109453
$ws.Range("O18").Value = 0.06965174129353234
$ws.Range("S18").Value = 0.1343283582089552
$ws.Range("F19").Value = 0.01850424055512722
$ws.Range("H19").Value = 0.2035466461063994
$ws.Range("I19").Value = 0.09252120277563608
$ws.Range("J19").Value = 0.3646877409406322
$ws.Range("K19").Value = 0.1002313030069391
$ws.Range("M19").Value = 0.01927525057825752
$ws.Range("N19").Value = 0.0007710100231303007
$ws.Range("O19").Value = 0.07324595219737856
$ws.Range("S19").Value = 0.1272166538164996
